$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 87.7
$ws.Range("C4").Value = 7
$ws.Range("E12").ClearContents()
$ws.Range("B14").Value = 265400
